$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new placeholder values to row 2 (columns G and J)
# Single-quoted strings are used so that the literal '${...}' text is
# written to the cell instead of being treated as PowerShell variable
# expansion syntax.
$ws.Range("G2").Value = '${emprest}'
$ws.Range("J2").Value = '${quantity_requested}'

# Update the active selection to I9 as recorded in the workbook view
$ws.Range("I9").Select()
